$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell A1 also picks up the date number format style
$ws.Cells.Item(1, 1).NumberFormat = "m/d/yy h:mm"

# Update column A width to fit the new date values
$ws.Columns.Item(1).ColumnWidth = 14

# Row 2 data
$ws.Cells.Item(2, 1).Value = 42605.671006944445
$ws.Cells.Item(2, 1).NumberFormat = "m/d/yy h:mm"

$ws.Cells.Item(2, 2).Value = -16
$ws.Cells.Item(2, 3).Value = 64
$ws.Cells.Item(2, 4).Value = 34
$ws.Cells.Item(2, 5).Value = 11
$ws.Cells.Item(2, 6).Value = 88
$ws.Cells.Item(2, 7).Value = 2694
$ws.Cells.Item(2, 8).Value = 5056
$ws.Cells.Item(2, 9).Value = 596
$ws.Cells.Item(2, 10).Value = 92
$ws.Cells.Item(2, 11).Value = 49
$ws.Cells.Item(2, 12).Value = 2
$ws.Cells.Item(2, 13).Value = 15
$ws.Cells.Item(2, 14).Value = "Bag"
